# Commit: "Kayıt silindi: 1480" -> delete the record whose Kayıt No (column A)
# equals "1480" from the master "Kayitlar" sheet, and from the matching
# district sheet "Merkez İlçe" (the district recorded for that row), which
# mirrors a filtered subset of "Kayitlar".

$wb = $excel.ActiveWorkbook

function Remove-KayitRow($SheetName, $KayitNo) {
    $ws = $wb.Worksheets.Item($SheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cellValue = $ws.Cells.Item($r, 1).Text
        if ($cellValue -eq $KayitNo) {
            $ws.Rows($r).Delete()
            Write-Output "Deleted row $r (Kayıt No $KayitNo) from '$SheetName'"
            return $true
        }
    }

    Write-Output "Kayıt No $KayitNo not found in '$SheetName'"
    return $false
}

Remove-KayitRow "Kayitlar" "1480"
Remove-KayitRow "Merkez İlçe" "1480"
